{"js": "async (context) => {\n  const replacements = [\n    [\"2025-11-12 Wednesday\", \"2025-11-13 Thursday\"],\n    [\"663\u00f78=\", \"280\u00f76=\"],\n    [\"669\u00f78=\", \"996\u00f74=\"],\n    [\"409\u00f78=\", \"532\u00f72=\"],\n    [\"559\u00f72=\", \"753\u00f77=\"],\n    [\"347\u00f74=\", \"930\u00f73=\"],\n    [\"148\u00f72=\", \"247\u00f77=\"],\n    [\"550\u00f79=\", \"437\u00f74=\"],\n    [\"710\u00f76=\", \"966\u00f72=\"],\n    [\"362\u00f74=\", \"746\u00f79=\"],\n    [\"885\u00f73=\", \"902\u00f75=\"],\n    [\"914\u00f74=\", \"976\u00f74=\"],\n    [\"470\u00f76=\", \"695\u00f73=\"],\n    [\"356\u00f73=\", \"156\u00f79=\"],\n    [\"939\u00f73=\", \"545\u00f76=\"],\n    [\"830\u00f72=\", \"494\u00f79=\"],\n    [\"303\u00f75=\", \"104\u00f73=\"],\n    [\"482\u00f73=\", \"942\u00f79=\"],\n    [\"324\u00f74=\", \"886\u00f72=\"],\n    [\"574\u00f78=\", \"474\u00f76=\"],\n    [\"316\u00f76=\", \"717\u00f75=\"],\n    [\"812\u00f78=\", \"306\u00f74=\"],\n    [\"405\u00f73=\", \"196\u00f77=\"],\n    [\"102\u00f79=\", \"561\u00f76=\"],\n    [\"687\u00f75=\", \"113\u00f75=\"],\n    [\"879\u00f73=\", \"560\u00f76=\"],\n  ];\n\n  const body = context.document.body;\n\n  for (const [oldText, newText] of replacements) {\n    const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n    results.load(\"items\");\n    await context.sync();\n\n    for (const range of results.items) {\n      range.insertText(newText, \"Replace\");\n    }\n    await context.sync();\n  }\n};\n", "ps1": "$d = $word.ActiveDocument\n\n$pairs = @(\n    @(\"2025-11-12 Wednesday\", \"2025-11-13 Thursday\"),\n    @(\"663\u00f78=\", \"280\u00f76=\"),\n    @(\"669\u00f78=\", \"996\u00f74=\"),\n    @(\"409\u00f78=\", \"532\u00f72=\"),\n    @(\"559\u00f72=\", \"753\u00f77=\"),\n    @(\"347\u00f74=\", \"930\u00f73=\"),\n    @(\"148\u00f72=\", \"247\u00f77=\"),\n    @(\"550\u00f79=\", \"437\u00f74=\"),\n    @(\"710\u00f76=\", \"966\u00f72=\"),\n    @(\"362\u00f74=\", \"746\u00f79=\"),\n    @(\"885\u00f73=\", \"902\u00f75=\"),\n    @(\"914\u00f74=\", \"976\u00f74=\"),\n    @(\"470\u00f76=\", \"695\u00f73=\"),\n    @(\"356\u00f73=\", \"156\u00f79=\"),\n    @(\"939\u00f73=\", \"545\u00f76=\"),\n    @(\"830\u00f72=\", \"494\u00f79=\"),\n    @(\"303\u00f75=\", \"104\u00f73=\"),\n    @(\"482\u00f73=\", \"942\u00f79=\"),\n    @(\"324\u00f74=\", \"886\u00f72=\"),\n    @(\"574\u00f78=\", \"474\u00f76=\"),\n    @(\"316\u00f76=\", \"717\u00f75=\"),\n    @(\"812\u00f78=\", \"306\u00f74=\"),\n    @(\"405\u00f73=\", \"196\u00f77=\"),\n    @(\"102\u00f79=\", \"561\u00f76=\"),\n    @(\"687\u00f75=\", \"113\u00f75=\"),\n    @(\"879\u00f73=\", \"560\u00f76=\")\n)\n\nforeach ($pair in $pairs) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n    $rng = $d.Content\n    $rng.Find.Execute($oldText, $false, $false, $false, $false, $false, $true, 1, $false, $newText, 2) | Out-Null\n}\n"}
